$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Named Characters" - just a "Wargear" header in the new D column
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("D1").Value = "Wargear"
$ws1.Range("D1").Select()

# ---------------------------------------------------------------------------
# Sheet 2: "HQ"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("D1").Value = "Wargear"
$ws2.Range("D2").Value = "Gauss cannon, Staff of light"
$ws2.Range("D3").Value = "Staff of light"
$ws2.Range("D4").Value = "Staff of light"
$ws2.Range("D5").Value = "Staff of light"
$ws2.Range("D6").Value = "Staff of light"
$ws2.Columns.Item(2).ColumnWidth = 14.5
$ws2.Columns.Item(3).ColumnWidth = 15.5
$ws2.Range("D1").Select()

# ---------------------------------------------------------------------------
# Sheet 3: "Troops"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("D1").Value = "Wargear"
$ws3.Range("D2").Value = "Gauss blaster"
$ws3.Range("D3").Value = "Gauss flayer"
$ws3.Range("D3").Select()

# ---------------------------------------------------------------------------
# Sheet 4: "Elites"
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("D1").Value = "Wargear"
$ws4.Range("D2").Value = "Synaptic disintergrator"
$ws4.Range("D3").Value = "Flayer claws"
$ws4.Range("D4").Value = "Rod of covenant"
$ws4.Range("D5").Value = "Warscythe"
$ws4.Range("D6").Value = "Heat ray, Massive forelimbs"
$ws4.Range("D6").Select()

# ---------------------------------------------------------------------------
# Sheet 5: "Fast Attack"
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("D1").Value = "Wargear"
$ws5.Range("D2").Value = "Feeder mandibles"
$ws5.Range("D3").Value = "Viscious claws"
$ws5.Range("D4").Value = "Gauss cannon"
$ws5.Range("D5").Value = "2 Gauss blasters"
$ws5.Columns.Item(3).ColumnWidth = 15.5
$ws5.Range("D5").Select()

# ---------------------------------------------------------------------------
# Sheet 6: "Heavy Support" (no header cell was added to D1 in this sheet)
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("D2").Value = "Gauss cannon, Twin tesla destructor"
$ws6.Range("D3").Value = "Automaton claws"
$ws6.Range("D4").Value = "2*Gauss flayer arrays, Doomsday cannon"
$ws6.Range("D5").Value = "Heavy gauss cannon"
$ws6.Range("D6").Value = "4*Gauss flux arcs, Particle whip"
$ws6.Range("D7").Value = "Crackling tendrils"
$ws6.Range("D4").Select()

# ---------------------------------------------------------------------------
# Sheet 8: "Flyers" (selected before sheet 7 so that sheet 7 ends up active)
# ---------------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item(8)
$ws8.Range("D1").Value = "Wargear"
$ws8.Range("D2").Value = "2*Tesla destructors, Death ray"
$ws8.Range("D3").Value = "2*Tesla destructors"
$ws8.Range("D1").Select()

# ---------------------------------------------------------------------------
# Sheet 7: "Dedicated Transports" (this is the sheet left active/selected)
# ---------------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item(7)
$ws7.Range("D1").Value = "Wargear"
$ws7.Range("D2").Value = "2*Gauss flayer arrays"
$ws7.Range("D1").Select()
